# Applies the commit's changes:
#  1. Insert a new "Player Info" sheet as the first sheet, with an
#     ID/NAME/BATTING_HAND/BOWL_STYLE header row (bold/centered/bordered,
#     matching the existing header style) and one data row for player 4400.
#  2. On "ODI Batting": rename column MATCH_CARD_LINK -> MATCH_CODE and
#     replace the scorecard URL with the bare match code (4485).
#  3. On "ODI Bowling": same MATCH_CARD_LINK -> MATCH_CODE rename/value swap.

$wb = $excel.ActiveWorkbook

# --- ODI Batting: MATCH_CARD_LINK -> MATCH_CODE -----------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").NumberFormat = "@"
$batting.Range("D2").Value = "4485"
$batting.Range("D2").Style = "Normal"

# --- ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE -----------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").NumberFormat = "@"
$bowling.Range("B2").Value = "4485"
$bowling.Range("B2").Style = "Normal"

# --- New "Player Info" sheet, placed before "ODI Batting" -------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4400"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "Nitish Rana"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# Match the bold / centered / thin-bordered header style already used on the
# other two sheets.
$header = $playerInfo.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108  # xlCenter
$header.VerticalAlignment = -4160    # xlTop
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2
